$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Style update: the title (row 1) and the table header (row 2) end up
#    sharing the SAME font: bold + white text. The title drops its old
#    explicit 14pt size (reverting to the workbook default size) and the
#    header keeps its size but becomes white, so the two converge onto a
#    single shared font entry instead of two separate ones.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # "Training Dashboard"
$ws2 = $wb.Worksheets.Item(2)   # "Exam Dashboard"

$white = 16777215   # RGB(255,255,255) -> 0xFFFFFF, Excel's BGR-packed long

# Normalize the title's font size first so it converges with the header's
# font before we recolor both -- this keeps the font table from growing
# more than necessary.
$ws1.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Size = 11

# Recolor the title and the header row on both sheets to bold white text
# (fill/border stay untouched).
$ws1.Range("A1").Font.Color = $white
$ws1.Range("A2:K2").Font.Color = $white

$ws2.Range("A1").Font.Color = $white
$ws2.Range("A2:G2").Font.Color = $white

# ---------------------------------------------------------------------
# 2) Data refresh on the Training Dashboard sheet: the "PERIOD TO EXPIRE"
#    and "LAST UPDATE" columns move forward for the first two trainings.
#    The leading apostrophe forces the "LAST UPDATE" cells to stay plain
#    text (e.g. "16-Sep-2025") instead of being auto-parsed into a date
#    serial, matching how the rest of the date-like columns in this sheet
#    are stored.
# ---------------------------------------------------------------------
$ws1.Range("H3").Value = -29
$ws1.Range("I3").Value = "'16-Sep-2025"

$ws1.Range("H4").Value = -41
$ws1.Range("I4").Value = "'16-Sep-2025"

$wb.Save()
